$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.654.22'
$ws.Range('E2').Value = '  +1.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.891.26'
$ws.Range('E3').Value = '  +1.79%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.22'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4828'
$ws.Range('E7').Value = '  +0.92%  '
$ws.Range('E8').Value = '  +2.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06556'
$ws.Range('E9').Value = '  +1.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.874.31'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.80'
$ws.Range('E11').Value = '  +3.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07451'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.097'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.78'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6669'
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.608.99'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.22'
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007575'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '231.31'
$ws.Range('E20').Value = '  +3.70%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.104.64'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.269'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.192'
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.382'
$ws.Range('E25').Value = '  +1.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.38'
$ws.Range('E26').Value = '  +2.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.70'
$ws.Range('E27').Value = '  +1.53%  '
$ws.Range('E28').Value = '  +1.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1021'
$ws.Range('E29').Value = '  +11.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.396'
$ws.Range('E30').Value = '  -2.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.328'
$ws.Range('E31').Value = '  +2.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.025'
$ws.Range('E32').Value = '  +1.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05051'
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.206'
$ws.Range('E34').Value = '  +5.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7498'
$ws.Range('E35').Value = '  +3.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.708'
$ws.Range('E37').Value = '  +0.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01883'
$ws.Range('E38').Value = '  +2.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.651'
$ws.Range('E39').Value = '  +1.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9212'
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('E41').Value = '  +1.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.90'
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4289'
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.650'
$ws.Range('E45').Value = '  -3.87%  '
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1278'
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('E49').Value = '  -0.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.946'
$ws.Range('E50').Value = '  +3.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.03'
$ws.Range('E51').Value = '  +1.03%  '
